$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values are stored as TEXT in the original workbook even
# when they look numeric (e.g. "0.999", "577.12"). Force text storage by
# setting the cell number format to Text ("@") before assigning the value,
# then restore the default "Normal" style so no stray style index is left
# behind on the cell (matching the un-styled cells in the source file).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "61.821.66"
$ws.Range("E2").Value = "  -1.60%  "
Set-TextValue "D3" "2.443.37"
$ws.Range("E3").Value = "  -0.20%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "577.12"
$ws.Range("E5").Value = "  -0.11%  "
Set-TextValue "D6" "141.22"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.21%  "
Set-TextValue "D9" "2.437.96"
$ws.Range("E9").Value = "  -0.26%  "
Set-TextValue "D10" "0.108"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +2.52%  "
Set-TextValue "D12" "5.16"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("E13").Value = "  -3.41%  "
Set-TextValue "D14" "25.84"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("E15").Value = "  +0.63%  "
Set-TextValue "D16" "0.0000171"
$ws.Range("E16").Value = "  -2.11%  "
Set-TextValue "D17" "61.794.08"
$ws.Range("E17").Value = "  -1.36%  "
Set-TextValue "D18" "2.439.26"
$ws.Range("E18").Value = "  -0.25%  "
Set-TextValue "D19" "10.57"
$ws.Range("E19").Value = "  -5.01%  "
Set-TextValue "D20" "7.18"
$ws.Range("E20").Value = "  +0.26%  "
Set-TextValue "D21" "324.09"
$ws.Range("E21").Value = "  -1.90%  "
Set-TextValue "D22" "4.05"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("E23").Value = "  -0.13%  "
Set-TextValue "D24" "1.90"
$ws.Range("E24").Value = "  -6.77%  "
Set-TextValue "D25" "64.99"
$ws.Range("E25").Value = "  -1.16%  "
Set-TextValue "D26" "9.17"
$ws.Range("E26").Value = "  +0.91%  "
Set-TextValue "D27" "583.23"
$ws.Range("E27").Value = "  -9.07%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D28" "2.545.80"
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D29" "1.00"
$ws.Range("E29").Value = "  +0.17%  "
Set-TextValue "D30" "0.0₃0929"
$ws.Range("E30").Value = "  -4.49%  "
Set-TextValue "D31" "7.84"
$ws.Range("E31").Value = "  -3.28%  "
Set-TextValue "D32" "1.37"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("E35").Value = "  +0.16%  "
Set-TextValue "D36" "4.74"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "152.79"
$ws.Range("E37").Value = "  +3.81%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D38" "0.373"
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("E39").Value = "  -5.28%  "
Set-TextValue "D40" "18.27"
$ws.Range("E40").Value = "  -1.61%  "
Set-TextValue "D41" "5.13"
$ws.Range("E41").Value = "  -2.80%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D42" "0.999"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D43" "42.13"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  -4.85%  "
Set-TextValue "D45" "2.36"
$ws.Range("E45").Value = "  -6.03%  "
Set-TextValue "D46" "0.0₆0276"
$ws.Range("E46").Value = "  +16.29%  "
Set-TextValue "D47" "140.64"
$ws.Range("E47").Value = "  -3.30%  "
Set-TextValue "D48" "3.55"
$ws.Range("E48").Value = "  -4.51%  "
Set-TextValue "D49" "0.597"
$ws.Range("E49").Value = "  -0.91%  "
Set-TextValue "D50" "0.0508"
$ws.Range("E50").Value = "  -3.40%  "
Set-TextValue "D51" "19.47"
$ws.Range("E51").Value = "  -1.93%  "
